$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Controls screen feature is now implemented ("Yes"), with a Location note.
$ws.Range("B26").Value = "Yes"
$ws.Range("C26").Value = "GameStateManager, Resources Menus, Resource Level archetypes Menus. Accessible ingame via main menu and pause menu."

# Update the view state: zoom level and current selection.
$excel.ActiveWindow.Zoom = 130
$ws.Range("C9").Select()
